# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) figures
# to match the latest data pull (GitHub Actions symbol-list update).
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the source data, which stores numbers/percentages as strings)
# instead of re-interpreting them as numeric/percentage values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'290.49"
$ws.Range("E2").Value = "'-0.04%"

# Row 3: OKB
$ws.Range("D3").Value = "'31.04"
$ws.Range("E3").Value = "'0.67%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'4.919"
$ws.Range("E4").Value = "'-0.36%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.07395"
$ws.Range("E5").Value = "'2.67%"

# Row 6: FTXToken
$ws.Range("D6").Value = "'2.208"
$ws.Range("E6").Value = "'23.59%"

# Row 7: KuCoinToken
$ws.Range("D7").Value = "'7.724"
$ws.Range("E7").Value = "'0.75%"

# Row 8: GateToken
$ws.Range("D8").Value = "'3.752"
$ws.Range("E8").Value = "'0.10%"

# Row 9: MXToken
$ws.Range("D9").Value = "'0.9125"
$ws.Range("E9").Value = "'1.92%"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.08906"
$ws.Range("E10").Value = "'15.09%"

# Row 11: WazirX
$ws.Range("E11").Value = "'2.22%"

# Row 12: MandalaExchangeToken
$ws.Range("D12").Value = "'0.08242"
$ws.Range("E12").Value = "'2.67%"

# Row 13: BitrueCoin
$ws.Range("D13").Value = "'0.03126"
$ws.Range("E13").Value = "'1.88%"

# Row 14: BitMartToken
$ws.Range("D14").Value = "'0.09982"
$ws.Range("E14").Value = "'-0.34%"

# Row 15: BitForexToken
$ws.Range("D15").Value = "'0.001502"
$ws.Range("E15").Value = "'-0.22%"

# Row 16: TigerCash
$ws.Range("D16").Value = "'0.005863"
$ws.Range("E16").Value = "'3.07%"

# Row 17: LEO
$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'0.61%"

# Row 18: BTSEToken
$ws.Range("D18").Value = "'2.086"
$ws.Range("E18").Value = "'0.15%"

# Row 19: BitpandaEcosystemToken
$ws.Range("D19").Value = "'0.3326"
$ws.Range("E19").Value = "'1.45%"

# Row 20: ProBitToken
$ws.Range("E20").Value = "'0.02%"

# Row 21: MCDex
$ws.Range("D21").Value = "'3.974"
$ws.Range("E21").Value = "'-1.89%"

# Row 22: ZBToken
$ws.Range("D22").Value = "'0.2189"
$ws.Range("E22").Value = "'4.25%"

# Row 23: CoinExToken
$ws.Range("D23").Value = "'0.04565"
$ws.Range("E23").Value = "'1.13%"

# Row 24: BitKan
$ws.Range("E24").Value = "'0.16%"

# Row 25: HotbitToken
$ws.Range("D25").Value = "'0.004582"
$ws.Range("E25").Value = "'14.29%"

# Row 26: NitroEx
$ws.Range("D26").Value = "'0.0001303"
$ws.Range("E26").Value = "'4.14%"

# Row 27: UpBots
$ws.Range("D27").Value = "'0.0003401"

# Row 39: One
$ws.Range("D39").Value = "'0.01587"
$ws.Range("E39").Value = "'-0.71%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.04467"
$ws.Range("E40").Value = "'1.93%"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.007343"
$ws.Range("E41").Value = "'-0.18%"

# Row 42: Dexo
$ws.Range("D42").Value = "'0.009558"
$ws.Range("E42").Value = "'24.96%"

# Row 43: BKEXToken
$ws.Range("D43").Value = "'0.1325"
$ws.Range("E43").Value = "'1.32%"

# Row 44: CEJI
$ws.Range("D44").Value = "'0.002315"
$ws.Range("E44").Value = "'12.29%"

# Row 45: LocalTraders
$ws.Range("D45").Value = "'0.008229"
$ws.Range("E45").Value = "'-10.65%"

# Row 46: CoinLion
$ws.Range("D46").Value = "'0.00006097"
$ws.Range("E46").Value = "'2.90%"

# Row 47: Kangarootoken
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.19%"

# Row 48: BOLO
$ws.Range("D48").Value = "'2.203"
$ws.Range("E48").Value = "'-1.89%"

# Row 49: CoinbaseStockToken
$ws.Range("D49").Value = "'0.002004"

# Row 50: CryptobidCoin
$ws.Range("E50").Value = "'0.19%"

# Row 51: SpecialPowerGold
$ws.Range("E51").Value = "'0.19%"
